$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.01432564932107159, 6.991465639761399, 1.110493442951626)
    3  = @(0.01304595685849999, 46.29620382419033, 0.890130492539376)
    4  = @(0.006174628445307307, 76.92516958383362, 1.451025732973183)
    5  = @(0.01005401168578613, 116.4452259784817, 1.237896088180598)
    6  = @(0.007853472290746061, 147.7555885713755, 0.6501546188889854)
    7  = @(0.007301446692124998, 187.3139259642583, 0.8203957453608403)
    8  = @(0.01341278535509757, 220.1113820478061, 1.190100339035406)
    9  = @(0.01043981904389949, 258.0971497528632, 0.8843296927132026)
    10 = @(0.01490267155112034, 287.7071415389927, 0.6127799019760743)
    11 = @(0.005572718299008441, 327.0803682114677, 1.417661879008176)
    12 = @(0.005978756055694598, 359.7540018323455, 1.086365417222171)
    13 = @(0.01146830785530686, 395.5575746271683, 0.7528285262379067)
    14 = @(0.006650313773219603, 430.9791554073337, 1.041036089718862)
    15 = @(0.008498604919453815, 466.9197757920596, 0.5175774094183757)
    16 = @(0.01214933101629917, 499.6142082516855, 1.0401307229293)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
